$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string slot reuse: process FilesTab (B4) first,
# then SamplesTab (B3), then CasesTab (B2) last, so the rebuilt shared
# string table lands the new query text at the same indices Excel produced.

# FilesTab query (B4): append ORDER BY / LIMIT clause
$filesQuery = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = $filesQuery + "`n order By f.file_name ASC LIMIT 100"

# SamplesTab query (B3): append ORDER BY / LIMIT clause
$samplesQuery = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"

# CasesTab query (B2): append ORDER BY / LIMIT clause
$casesQuery = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100"
